# Updated cryptos list on Sat Mar 16 14:23:42 UTC 2024 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures for the
# crypto listing on the active sheet, and swaps the FirstDigitalUSD / TheGraph
# rows (38 and 39) to reflect their new rank/order along with their own
# refreshed data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every cell this edit touches, with its new text value. All of these are
# stored as plain text in the sheet (prices like "1.00" / "0.999", and
# padded percentages like "  -0.31%  "), so we force a text number format
# before writing so Excel doesn't silently coerce them into numbers (which
# would drop significant trailing zeros, turn "1.00" into "1", etc.), then
# restore the default "Normal" style afterwards so no stray formatting is
# left behind.
$updates = @(
    @{Cell='D2'; Value='68.003.90'},
    @{Cell='E2'; Value='  -0.05%  '},
    @{Cell='D3'; Value='3.662.33'},
    @{Cell='E3'; Value='  -0.97%  '},
    @{Cell='D4'; Value='0.999'},
    @{Cell='E4'; Value='  -0.18%  '},
    @{Cell='D5'; Value='598.08'},
    @{Cell='E5'; Value='  +2.17%  '},
    @{Cell='D6'; Value='190.52'},
    @{Cell='E6'; Value='  +5.03%  '},
    @{Cell='D7'; Value='0.620'},
    @{Cell='E7'; Value='  -1.14%  '},
    @{Cell='E8'; Value='  +0.31%  '},
    @{Cell='D9'; Value='0.701'},
    @{Cell='E9'; Value='  -1.52%  '},
    @{Cell='D10'; Value='57.41'},
    @{Cell='E10'; Value='  +6.32%  '},
    @{Cell='D11'; Value='0.153'},
    @{Cell='E11'; Value='  -5.25%  '},
    @{Cell='D12'; Value='0.0000275'},
    @{Cell='E12'; Value='  -5.46%  '},
    @{Cell='D13'; Value='10.19'},
    @{Cell='E13'; Value='  -2.13%  '},
    @{Cell='D14'; Value='4.241.66'},
    @{Cell='E14'; Value='  -1.16%  '},
    @{Cell='D15'; Value='3.654.26'},
    @{Cell='E15'; Value='  -1.35%  '},
    @{Cell='E16'; Value='  +0.82%  '},
    @{Cell='D17'; Value='18.88'},
    @{Cell='E17'; Value='  -2.55%  '},
    @{Cell='D18'; Value='1.11'},
    @{Cell='E18'; Value='  -0.84%  '},
    @{Cell='D19'; Value='67.747.13'},
    @{Cell='E19'; Value='  -0.03%  '},
    @{Cell='D20'; Value='12.45'},
    @{Cell='E20'; Value='  -2.74%  '},
    @{Cell='D21'; Value='400.45'},
    @{Cell='E21'; Value='  -1.56%  '},
    @{Cell='D22'; Value='4.42'},
    @{Cell='E22'; Value='  -1.46%  '},
    @{Cell='D23'; Value='87.39'},
    @{Cell='E23'; Value='  -1.02%  '},
    @{Cell='D24'; Value='11.27'},
    @{Cell='E24'; Value='  +2.20%  '},
    @{Cell='D25'; Value='2.95'},
    @{Cell='E25'; Value='  -2.99%  '},
    @{Cell='D26'; Value='12.46'},
    @{Cell='E26'; Value='  -2.19%  '},
    @{Cell='D27'; Value='6.05'},
    @{Cell='E27'; Value='  -0.11%  '},
    @{Cell='D28'; Value='3.66'},
    @{Cell='E28'; Value='  -5.47%  '},
    @{Cell='D29'; Value='9.27'},
    @{Cell='E29'; Value='  -2.40%  '},
    @{Cell='D30'; Value='31.88'},
    @{Cell='E30'; Value='  -1.98%  '},
    @{Cell='D31'; Value='7.38'},
    @{Cell='E31'; Value='  -2.10%  '},
    @{Cell='D32'; Value='12.31'},
    @{Cell='E32'; Value='  -1.26%  '},
    @{Cell='D33'; Value='45.05'},
    @{Cell='E33'; Value='  +4.41%  '},
    @{Cell='D34'; Value='66.39'},
    @{Cell='E34'; Value='  +1.83%  '},
    @{Cell='D35'; Value='0.117'},
    @{Cell='E35'; Value='  +0.44%  '},
    @{Cell='D36'; Value='611.64'},
    @{Cell='E36'; Value='  +1.19%  '},
    @{Cell='D37'; Value='1.00'},
    @{Cell='E37'; Value='  +0.16%  '},
    @{Cell='B38'; Value='TheGraph'},
    @{Cell='C38'; Value='https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'},
    @{Cell='D38'; Value='0.393'},
    @{Cell='E38'; Value='  -1.10%  '},
    @{Cell='B39'; Value='FirstDigitalUSD'},
    @{Cell='C39'; Value='https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'},
    @{Cell='D39'; Value='0.999'},
    @{Cell='E39'; Value='  -0.21%  '},
    @{Cell='D40'; Value='0.0₃0777'},
    @{Cell='E40'; Value='  -13.05%  '},
    @{Cell='D41'; Value='0.135'},
    @{Cell='E41'; Value='  -0.78%  '},
    @{Cell='D42'; Value='2.90'},
    @{Cell='E42'; Value='  -3.04%  '},
    @{Cell='D43'; Value='0.0426'},
    @{Cell='E43'; Value='  -1.79%  '},
    @{Cell='D44'; Value='2.56'},
    @{Cell='E44'; Value='  -8.48%  '},
    @{Cell='D45'; Value='0.136'},
    @{Cell='E45'; Value='  +2.02%  '},
    @{Cell='D46'; Value='2.792.61'},
    @{Cell='E46'; Value='  -0.63%  '},
    @{Cell='E47'; Value='  +2.13%  '},
    @{Cell='D48'; Value='8.83'},
    @{Cell='E48'; Value='  -4.10%  '},
    @{Cell='D49'; Value='143.24'},
    @{Cell='E49'; Value='  +3.25%  '},
    @{Cell='D50'; Value='2.61'},
    @{Cell='E50'; Value='  -2.83%  '},
    @{Cell='D51'; Value='2.48'},
    @{Cell='E51'; Value='  -16.50%  '}
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.Style = "Normal"
}
